$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data in columns D (Price) and E (Volume(1h)) is plain text
# (e.g. "58.941.93", "  -0.46%  ") rather than numeric values, so force the
# text number format on the range before writing the new values. This stops
# Excel's automatic "looks like a number" conversion from mangling values
# such as "59.147.07" or rounding "137.95" to a binary float. The original
# style is restored afterwards so cell formatting is left untouched.
$dataRange = $ws.Range("D2:E51")
$origStyle = $dataRange.Style
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '59.147.07'
$ws.Range("E2").Value = '  -0.27%  '
$ws.Range("D3").Value = '2.522.19'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '536.12'
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").Value = '137.95'
$ws.Range("E6").Value = '  -1.11%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  +0.78%  '
$ws.Range("D9").Value = '2.521.14'
$ws.Range("E9").Value = '  -0.29%  '
$ws.Range("D10").Value = '0.101'
$ws.Range("E10").Value = '  +0.55%  '
$ws.Range("D11").Value = '0.158'
$ws.Range("E11").Value = '  -2.04%  '
$ws.Range("E12").Value = '  -1.81%  '
$ws.Range("E13").Value = '  -2.06%  '
$ws.Range("D14").Value = '2.971.80'
$ws.Range("E14").Value = '  +0.10%  '
$ws.Range("D15").Value = '23.09'
$ws.Range("E15").Value = '  -1.42%  '
$ws.Range("D16").Value = '59.071.20'
$ws.Range("E16").Value = '  -0.25%  '
$ws.Range("D17").Value = '0.0000139'
$ws.Range("E17").Value = '  -1.29%  '
$ws.Range("D18").Value = '2.523.56'
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("E19").Value = '  +0.46%  '
$ws.Range("D20").Value = '4.30'
$ws.Range("E20").Value = '  -0.57%  '
$ws.Range("D21").Value = '325.72'
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = '5.96'
$ws.Range("E23").Value = '  +2.24%  '
$ws.Range("E24").Value = '  +4.92%  '
$ws.Range("D25").Value = '0.424'
$ws.Range("E25").Value = '  -1.12%  '
$ws.Range("D26").Value = '0.167'
$ws.Range("E26").Value = '  +0.69%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").Value = '7.66'
$ws.Range("E28").Value = '  -2.29%  '
$ws.Range("D29").Value = '6.72'
$ws.Range("E29").Value = '  -1.89%  '
$ws.Range("D30").Value = '0.0₃0773'
$ws.Range("E30").Value = '  -0.43%  '
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("D32").Value = '1.18'
$ws.Range("E32").Value = '  +5.73%  '
$ws.Range("D33").Value = '162.37'
$ws.Range("E33").Value = '  -1.69%  '
$ws.Range("D34").Value = '1.48'
$ws.Range("E34").Value = '  +0.55%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").Value = '18.49'
$ws.Range("E36").Value = '  -0.26%  '
$ws.Range("D37").Value = '4.13'
$ws.Range("E37").Value = '  -3.51%  '
$ws.Range("E38").Value = '  -2.19%  '
$ws.Range("D39").Value = '36.63'
$ws.Range("E39").Value = '  -0.81%  '
$ws.Range("D40").Value = '0.818'
$ws.Range("E40").Value = '  +0.55%  '
$ws.Range("D41").Value = '3.64'
$ws.Range("E41").Value = '  -1.06%  '
$ws.Range("D42").Value = '287.17'
$ws.Range("E42").Value = '  +2.31%  '
$ws.Range("D43").Value = '5.21'
$ws.Range("E43").Value = '  -1.58%  '
$ws.Range("D44").Value = '132.43'
$ws.Range("E44").Value = '  +8.09%  '
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.22%  '
$ws.Range("D46").Value = '0.610'
$ws.Range("E46").Value = '  +1.92%  '
$ws.Range("D47").Value = '10.89'
$ws.Range("D48").Value = '0.0932'
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("D49").Value = '0.0510'
$ws.Range("E49").Value = '  -0.67%  '
$ws.Range("E50").Value = '  -1.14%  '
$ws.Range("D51").Value = '17.40'
$ws.Range("E51").Value = '  -2.20%  '

$dataRange.Style = $origStyle
